# Replace the hard-coded "Creative Planning Legal" signature block with
# blank underscored signature lines (name / address line 1 / address line 2).

$d = $word.ActiveDocument

$d.Content.Find.Execute("Creative Planning Legal", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "_________________________________", 2)

$d.Content.Find.Execute("5454 W. 110th Street", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "_________________________________", 2)

$d.Content.Find.Execute("Overland Park, KS  66211", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "_________________________________", 2)
